# Two more ManagePage tests:
#   ChangePasswordwithInvalidNewPassword
#   ChangePasswordwithoutConfirmedNewPassword
#
# Adds two data rows (row 3 and row 4) to the existing test-data sheet,
# mirroring the shape of the already-present row 2, plus the associated
# mailto hyperlinks on the Username column, and updates the sheet
# selection / column-A width to reflect the wider content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: ChangePasswordwithInvalidNewPassword -------------------------
$ws.Range("A3").Value = "ChangePasswordwithInvalidNewPassword"
$ws.Range("B3").Value = "stoyanski@mail.com"
$ws.Range("C3").Value = "password"
$ws.Range("D3").Value = " "
$ws.Range("E3").Value = "dada"
$ws.Range("F3").Value = "password"

# --- Row 4: ChangePasswordwithoutConfirmedNewPassword ---------------------
$ws.Range("A4").Value = "ChangePasswordwithoutConfirmedNewPassword"
$ws.Range("B4").Value = "stoyanski@mail.com"
$ws.Range("C4").Value = "password"
$ws.Range("D4").Value = "pass"
$ws.Range("E4").Value = " "
$ws.Range("F4").Value = "password"

# --- Hyperlinks on the new Username cells (same mailto as row 2) ---------
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:stoyanski@mail.com")
$ws.Range("B3").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:stoyanski@mail.com")
$ws.Range("B4").Style = "Hyperlink"

# --- Column A is now a bit wider because of the new, longer keys ---------
$ws.Columns.Item(1).ColumnWidth = 45.5

# --- Selection left where the author last clicked -------------------------
$ws.Range("E6").Select()
